$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "27.467.86"
Set-TextValue "E2" "  -3.69%  "
Set-TextValue "D3" "1.755.51"
Set-TextValue "E3" "  -2.90%  "
Set-TextValue "D4" "1.003"
Set-TextValue "E4" "  +0.28%  "
Set-TextValue "D5" "322.97"
Set-TextValue "E5" "  -1.76%  "
Set-TextValue "D6" "1.000"
Set-TextValue "E6" "  +0.10%  "
Set-TextValue "D7" "0.4421"
Set-TextValue "E7" "  -1.88%  "
Set-TextValue "D8" "0.3702"
Set-TextValue "E8" "  -1.80%  "
Set-TextValue "D9" "44.84"
Set-TextValue "E9" "  +0.14%  "
Set-TextValue "D10" "0.07713"
Set-TextValue "E10" "  +2.24%  "
Set-TextValue "D11" "1.110"
Set-TextValue "E11" "  -3.61%  "
Set-TextValue "D12" "1.000"
Set-TextValue "E12" "  -0.07%  "
Set-TextValue "D13" "21.52"
Set-TextValue "E13" "  -5.05%  "
Set-TextValue "D14" "6.145"
Set-TextValue "E14" "  -2.74%  "
Set-TextValue "D15" "7.401"
Set-TextValue "D16" "1.763.98"
Set-TextValue "E16" "  -2.17%  "
Set-TextValue "D17" "90.06"
Set-TextValue "E17" "  +11.39%  "
Set-TextValue "D18" "0.00001071"
Set-TextValue "E18" "  -1.96%  "
Set-TextValue "D19" "0.06233"
Set-TextValue "E19" "  -7.94%  "
Set-TextValue "D20" "1.000"
Set-TextValue "E20" "  +0.07%  "
Set-TextValue "D21" "17.32"
Set-TextValue "E21" "  -2.02%  "
Set-TextValue "D22" "6.163"
Set-TextValue "E22" "  -2.53%  "
Set-TextValue "D23" "0.5276"
Set-TextValue "E23" "  -3.06%  "
Set-TextValue "D24" "27.518.18"
Set-TextValue "E24" "  -3.42%  "
Set-TextValue "D25" "11.49"
Set-TextValue "E25" "  -2.95%  "
Set-TextValue "D26" "2.303"
Set-TextValue "E26" "  -4.37%  "
Set-TextValue "D27" "20.45"
Set-TextValue "E27" "  -0.95%  "
Set-TextValue "D28" "152.75"
Set-TextValue "E28" "  +0.80%  "
Set-TextValue "D29" "2.282"
Set-TextValue "E29" "  -3.53%  "
Set-TextValue "D30" "1.957.06"
Set-TextValue "E30" "  -2.47%  "
Set-TextValue "D31" "126.99"
Set-TextValue "E31" "  -4.59%  "
Set-TextValue "D32" "1.168"
Set-TextValue "E32" "  -7.91%  "
Set-TextValue "D33" "5.694"
Set-TextValue "E33" "  -2.54%  "
Set-TextValue "D34" "0.09180"
Set-TextValue "E34" "  -1.66%  "
Set-TextValue "D35" "3.626"
Set-TextValue "E35" "  -9.45%  "
Set-TextValue "D36" "12.56"
Set-TextValue "E36" "  +2.47%  "
Set-TextValue "D37" "0.02301"
Set-TextValue "E37" "  -1.12%  "
Set-TextValue "D38" "0.2149"
Set-TextValue "E38" "  -5.24%  "
Set-TextValue "D39" "0.06096"
Set-TextValue "E39" "  -4.20%  "
Set-TextValue "D42" "1.179"
Set-TextValue "E42" "  -2.79%  "
Set-TextValue "D43" "1.000"
Set-TextValue "E43" "  +0.09%  "
Set-TextValue "D44" "7.901"
Set-TextValue "E44" "  -2.91%  "
Set-TextValue "D45" "1.380"
Set-TextValue "E45" "  -4.71%  "
Set-TextValue "D46" "13.56"
Set-TextValue "E46" "  -2.01%  "
Set-TextValue "D47" "0.5940"
Set-TextValue "E47" "  -2.50%  "
Set-TextValue "D48" "3.710"
Set-TextValue "E48" "  -2.71%  "
Set-TextValue "D49" "125.85"
Set-TextValue "E49" "  -2.16%  "
Set-TextValue "D50" "1.974"
Set-TextValue "E50" "  -3.31%  "
Set-TextValue "D51" "0.06865"
Set-TextValue "E51" "  -3.28%  "

# Row 40/41 swap: coin identity + link swap, with updated price/volume
Set-TextValue "B40" "TheSandbox"
Set-TextValue "C40" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D40" "0.6426"
Set-TextValue "E40" "  -2.83%  "
Set-TextValue "B41" "InternetComputer(DFINITY)"
Set-TextValue "C41" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D41" "5.034"
Set-TextValue "E41" "  -2.66%  "
